$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.958
$ws.Range("C7").Value = -13.385
$ws.Range("B8").Value = 6.967000000000001
$ws.Range("B10").Value = 5.935
$ws.Range("E10").Value = 16.638
$ws.Range("B12").Value = 5.715
$ws.Range("E12").Value = 17.12
$ws.Range("E13").Value = 16.43
$ws.Range("E14").Value = 16.917
$ws.Range("C15").Value = -13.413
$ws.Range("B18").Value = 5.609
$ws.Range("C18").Value = -13.483
$ws.Range("D18").Value = -8.797000000000001
$ws.Range("D19").Value = -8.113000000000001
$ws.Range("C20").Value = -12.673
$ws.Range("D27").Value = -8.106
$ws.Range("C29").Value = -12.221
$ws.Range("E29").Value = 16.951
$ws.Range("C30").Value = -12.486
$ws.Range("C31").Value = -12.788
$ws.Range("D31").Value = -8.126999999999999
$ws.Range("E32").Value = 16.86
$ws.Range("E35").Value = 16.358
$ws.Range("B37").Value = 8.494
$ws.Range("D38").Value = -7.808
$ws.Range("C40").Value = -12.782
$ws.Range("D42").Value = -8.468
$ws.Range("E43").Value = 17.023
$ws.Range("D44").Value = -7.766
$ws.Range("D47").Value = -7.458999999999999
$ws.Range("E48").Value = 17.17599999999999
$ws.Range("E49").Value = 16.349
$ws.Range("C50").Value = -12.858
$ws.Range("E50").Value = 16.435
$ws.Range("B55").Value = 5.136
$ws.Range("E56").Value = 16.149
$ws.Range("D58").Value = -8.134
$ws.Range("D65").Value = -7.723999999999999
$ws.Range("B68").Value = 5.271000000000001
$ws.Range("C68").Value = -11.307
$ws.Range("E69").Value = 17.241
$ws.Range("D73").Value = -7.981000000000002
$ws.Range("C76").Value = -13.779
$ws.Range("B77").Value = 6.119
$ws.Range("B78").Value = 7.631
$ws.Range("B81").Value = 5.672
$ws.Range("E81").Value = 16.485
$ws.Range("B82").Value = 5.860000000000001
$ws.Range("C87").Value = -12.75
$ws.Range("C88").Value = -12.827
$ws.Range("D90").Value = -7.450999999999999
$ws.Range("E92").Value = 17.641
$ws.Range("D94").Value = -7.296000000000001
$ws.Range("D95").Value = -7.537000000000001
$ws.Range("C96").Value = -12.705
$ws.Range("C98").Value = -13.649
$ws.Range("C101").Value = -12.612
$ws.Range("D101").Value = -7.768000000000001
$ws.Range("C102").Value = -13.492
